$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new value in column C, row 11 (new shared string "Delivery Note")
$ws.Range("C11").Value = "Delivery Note"

# Update the active selection to C9
$ws.Range("C9").Select()
